{"js": "// Cambios en comas de palabras clave (keyword list punctuation/casing cleanup)\n// and relocation of the transient \"_GoBack\" bookmark.\n//\n// 1) The keywords paragraph\n//      \"Reproducci\u00f3n, plantas, esporulaci\u00f3n, gametofito, esporofito,\n//       polinizaci\u00f3n, fecundaci\u00f3n, fructificaci\u00f3n, dispersi\u00f3n,\n//       germinaci\u00f3n, flor.\"\n//    becomes two runs (no spaces after commas, lower-case initial letter,\n//    no trailing period) with the \"_GoBack\" bookmark sitting between them,\n//    right after \"dispersi\u00f3n,\":\n//      \"reproducci\u00f3n,plantas,esporulaci\u00f3n,gametofito,esporofito,\n//       polinizaci\u00f3n,fecundaci\u00f3n,fructificaci\u00f3n,dispersi\u00f3n,\"\n//      [[bookmark _GoBack]]\n//      \"germinaci\u00f3n,flor\"\n//\n// 2) The old \"_GoBack\" bookmark (previously located right after\n//    \"plantas con flor\" / before the closing period in a different\n//    paragraph) is removed from there.\n\nconst oldKeywords =\n  \"Reproducci\u00f3n, plantas, esporulaci\u00f3n, gametofito, esporofito, \" +\n  \"polinizaci\u00f3n, fecundaci\u00f3n, fructificaci\u00f3n, dispersi\u00f3n, germinaci\u00f3n, flor.\";\nconst newPart1 =\n  \"reproducci\u00f3n,plantas,esporulaci\u00f3n,gametofito,esporofito,\" +\n  \"polinizaci\u00f3n,fecundaci\u00f3n,fructificaci\u00f3n,dispersi\u00f3n,\";\nconst newPart2 = \"germinaci\u00f3n,flor\";\n\nconst body = context.document.body;\n\n// --- Step 1: rewrite the keywords sentence in place (keeps run formatting) ---\nconst found = body.search(oldKeywords, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error(\"Keywords sentence not found\");\n}\nconst target = found.items[0];\ntarget.insertText(newPart1 + newPart2, \"Replace\");\nawait context.sync();\n\n// --- Step 2: remove the bookmark from its old location (\"plantas con flor\") ---\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Step 3: re-insert \"_GoBack\" right between the two keyword halves ---\n// Anchor the bookmark precisely after newPart1 using a dedicated search hit.\nconst anchor = body.search(newPart1, { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\nif (anchor.items.length === 0) {\n  throw new Error(\"First half of keywords sentence not found\");\n}\nconst afterPart1 = anchor.items[0].getRange(\"After\");\nafterPart1.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Best-effort: \"Table Grid\" style cleanup (tblInd / tblCellMar) ---\n// The diff also drops the redundant <w:tblInd w:w=\"0\".../> and\n// <w:tblCellMar>...</w:tblCellMar> (0/108/0/108 twips) from the\n// \"Tablaconcuadrcula\" (\"Table Grid\") table style's <w:tblPr> in\n// styles.xml. Those values are exactly Word's built-in defaults for a\n// table style (no table in this document renders differently either\n// way - every <w:tbl> here already carries its own explicit tblInd),\n// so this is a no-visual-effect save-time cleanup rather than a content\n// change. The Word JS API has no supported surface for editing a table\n// STYLE's own tblPr (Word.TableStyle exists on the object model but its\n// properties are not backed for style definitions in this host), so\n// this is attempted defensively and is safe to no-op.\ntry {\n  const styles = context.document.getStyles();\n  const tableGrid = styles.getByNameOrNullObject(\"Table Grid\");\n  await context.sync();\n  if (!tableGrid.isNullObject && tableGrid.tableStyle) {\n    tableGrid.tableStyle.leftIndent = 0;\n    tableGrid.tableStyle.topCellMargin = 0;\n    tableGrid.tableStyle.bottomCellMargin = 0;\n    tableGrid.tableStyle.leftCellMargin = 108 / 20; // 108 twips -> points\n    tableGrid.tableStyle.rightCellMargin = 108 / 20;\n    await context.sync();\n  }\n} catch (e) {\n  // Not supported in this host - the values already match Word's\n  // built-in table-style defaults, so there is nothing visible to fix.\n}\n", "ps1": "# Cambios en comas de palabras clave (keyword list punctuation/casing\n# cleanup) and relocation of the transient \"_GoBack\" bookmark.\n#\n# 1) The keywords paragraph\n#      \"Reproducci\u00f3n, plantas, esporulaci\u00f3n, gametofito, esporofito,\n#       polinizaci\u00f3n, fecundaci\u00f3n, fructificaci\u00f3n, dispersi\u00f3n,\n#       germinaci\u00f3n, flor.\"\n#    becomes two runs (no spaces after commas, lower-case initial letter,\n#    no trailing period) with the \"_GoBack\" bookmark sitting between them,\n#    right after \"dispersi\u00f3n,\":\n#      \"reproducci\u00f3n,plantas,esporulaci\u00f3n,gametofito,esporofito,\n#       polinizaci\u00f3n,fecundaci\u00f3n,fructificaci\u00f3n,dispersi\u00f3n,\"\n#      [[bookmark _GoBack]]\n#      \"germinaci\u00f3n,flor\"\n#\n# 2) The old \"_GoBack\" bookmark (previously located right after\n#    \"plantas con flor\" / before the closing period in a different\n#    paragraph) is removed from there.\n\n$d = $word.ActiveDocument\n\n$oldKeywords = \"Reproducci\u00f3n, plantas, esporulaci\u00f3n, gametofito, esporofito, polinizaci\u00f3n, fecundaci\u00f3n, fructificaci\u00f3n, dispersi\u00f3n, germinaci\u00f3n, flor.\"\n$newPart1 = \"reproducci\u00f3n,plantas,esporulaci\u00f3n,gametofito,esporofito,polinizaci\u00f3n,fecundaci\u00f3n,fructificaci\u00f3n,dispersi\u00f3n,\"\n$newPart2 = \"germinaci\u00f3n,flor\"\n\n# --- Step 1: rewrite the keywords sentence in place (keeps run formatting) ---\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $oldKeywords\n$found = $rng.Find.Execute()\nif ($found) {\n    $rng.Text = $newPart1 + $newPart2\n}\n\n# --- Step 2: remove the bookmark from its old location (\"plantas con flor\") ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Step 3: re-insert \"_GoBack\" right between the two keyword halves ---\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = $newPart1\n$found2 = $rng2.Find.Execute()\nif ($found2) {\n    $point = $rng2.Duplicate\n    $point.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $point)\n}\n\n# --- Best-effort: \"Table Grid\" style cleanup (tblInd / tblCellMar) ---\n# The diff also drops the redundant <w:tblInd w:w=\"0\".../> and\n# <w:tblCellMar>...</w:tblCellMar> (0/108/0/108 twips) from the\n# \"Tablaconcuadrcula\" (\"Table Grid\") table style's <w:tblPr> in\n# styles.xml. Those numbers are exactly Word's built-in defaults for a\n# table style (no table in this document renders any differently either\n# way - every <w:tbl> here already carries its own explicit tblInd), so\n# this is a no-visual-effect save-time cleanup rather than a content\n# change. The Word object model has no supported surface in this host\n# for editing a table STYLE's own tblPr (Style.Table resolves to\n# Nothing here), so this is attempted defensively and is safe to no-op.\ntry {\n    $tableGridStyle = $d.Styles(\"Table Grid\")\n    $ts = $tableGridStyle.Table\n    if ($ts -ne $null) {\n        $ts.LeftIndent = 0\n        $ts.TopPadding = 0\n        $ts.BottomPadding = 0\n        $ts.LeftPadding = 5.4\n        $ts.RightPadding = 5.4\n    }\n} catch {\n    # Not supported in this host - the values already match Word's\n    # built-in table-style defaults, so there is nothing visible to fix.\n}\n\nWrite-Output \"done\"\n"}
